$d = $word.ActiveDocument

# 1. Replace the placeholder "your stuff" paragraph text with "Hello World!"
[void]$d.Content.Find.Execute(">>>  your stuff after this line >>>", $true, $false, $false, $false, $false, $true, 1, $false, "Hello World!", 2)

# 2. Remove the existing "_GoBack" bookmark (currently sitting in the title paragraph)
$existing = $d.Bookmarks("_GoBack")
$existing.Delete()

# 3. Re-create the "_GoBack" bookmark immediately after the new "Hello World!" text,
#    without splitting that run. We do this by temporarily appending a marker
#    character after the text (so the insertion point is no longer at the very
#    end of the paragraph, which the Range/Bookmark placement is unreliable for),
#    adding the bookmark right before the marker, and then removing the marker.
$p4 = $d.Paragraphs(4)
$rEnd = $p4.Range
[void]$rEnd.MoveEnd(1, -1)
$markerStart = $rEnd.End
[void]$rEnd.InsertAfter("Z")

$rb = $d.Range($markerStart, $markerStart)
[void]$d.Bookmarks.Add("_GoBack", $rb)

$rz = $d.Range($markerStart, $markerStart + 1)
[void]$rz.Delete()
